$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Formula = "=1.332267629550188E-15"
$ws.Range("B2").Value = $ws.Range("B2").Value2

$ws.Range("C2").Formula = "=2.509104035652854E-14"
$ws.Range("C2").Value = $ws.Range("C2").Value2

$ws.Range("D2").Value = 0.7127328510149897

$ws.Range("E2").Value = 198602002.3250627

$ws.Range("G2").Value = 198602003.0377955

# Row 3
$ws.Range("B3").Value = 0.02258322285507441

$ws.Range("C3").Value = 1.65323645889881

$ws.Range("D3").Value = 0.1529057820181812

$ws.Range("E3").Value = 6.48142807727062

$ws.Range("G3").Value = 8.310153541042686
